$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new header row at the top; existing data shifts down to rows 2-21
$ws.Rows.Item(1).Insert()

# Populate the new header row
$ws.Range("A1").Value = "Rank"
$ws.Range("B1").Value = "City Name"
$ws.Range("C1").Value = " Overnight International Visitor Spend (US`$ bn)  "
$ws.Range("D1").Value = "Year"

# Fill in the new Year column for every data row with 2013
$ws.Range("D2:D21").Value = 2013

# Reflect the new selection left behind after filling the Year column
$ws.Range("D2:D21").Select() | Out-Null
